# Ascender usuario a colaborador y casi todos los botones del admin funcionales
#
# Promote the "Nivel" (level) of several users in the Usuarios sheet:
#   - E3 (Id 2, user "qwe")                -> 1 to 2 (ascendido a colaborador)
#   - E4 (Id 3, José Manuel / Rex117)      -> 1 to 3 (ascendido a admin)
#   - E6 (Id 5, SebastianJerez)            -> 1 to 3 (ascendido a admin)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Usuarios")

$ws.Range("E3").Value = 2
$ws.Range("E4").Value = 3
$ws.Range("E6").Value = 3
